$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.947.47'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").Value = '1.672.63'
$ws.Range("E3").Value = '  +1.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.517'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.63%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").Value = '  +0.50%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.18'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.07%  '
$ws.Range("E11").Value = '  +1.16%  '
$ws.Range("D12").Value = '1.907.61'
$ws.Range("E12").Value = '  +1.02%  '
$ws.Range("D13").Value = '1.666.27'
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("E14").Value = '  -0.03%  '
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.57'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = '26.947.93'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '234.57'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.05%  '
$ws.Range("E19").Value = '  +3.70%  '
$ws.Range("E20").Value = '  -0.33%  '
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("E23").Value = '  -1.36%  '
$ws.Range("E24").Value = '  -2.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.81'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.40%  '
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.83%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  +0.22%  '
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("E32").Value = '  +0.77%  '
$ws.Range("D33").Value = '1.471.51'
$ws.Range("E33").Value = '  -4.91%  '
$ws.Range("E34").Value = '  +2.22%  '
$ws.Range("E35").Value = '  +1.39%  '
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("E37").Value = '  -0.59%  '
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("E39").Value = '  +1.11%  '
$ws.Range("E40").Value = '  +10.60%  '
$ws.Range("E41").Value = '  -3.54%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.29'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '66.78'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("D45").Value = '1.813.95'
$ws.Range("E45").Value = '  +1.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.780'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.61'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.64%  '
$ws.Range("E48").Value = '  +0.41%  '
$ws.Range("E49").Value = '  +2.10%  '
$ws.Range("E50").Value = '  +0.48%  '
$ws.Range("E51").Value = '  +0.41%  '
